$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1103
$ws.Cells.Item(1103, 1).Value = 44347
$ws.Cells.Item(1103, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1103, 2).Value = 'Reko193'
$ws.Cells.Item(1103, 3).Value = 3011
$ws.Cells.Item(1103, 4).Value = 'Reko Swish +447930169686'
$ws.Cells.Item(1103, 6).Value = 460.71

# Row 1104
$ws.Cells.Item(1104, 1).Value = 44347
$ws.Cells.Item(1104, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1104, 2).Value = 'Reko193'
$ws.Cells.Item(1104, 3).Value = 2611
$ws.Cells.Item(1104, 4).Value = 'Reko Swish +447930169686'
$ws.Cells.Item(1104, 6).Value = 55.29

# Row 1105
$ws.Cells.Item(1105, 1).Value = 44347
$ws.Cells.Item(1105, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1105, 2).Value = 'Reko193'
$ws.Cells.Item(1105, 3).Value = 1930
$ws.Cells.Item(1105, 4).Value = 'Reko Swish +447930169686'
$ws.Cells.Item(1105, 5).Value = 516

# Row 1106
$ws.Cells.Item(1106, 1).Value = 44347
$ws.Cells.Item(1106, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1106, 2).Value = 'Reko194'
$ws.Cells.Item(1106, 3).Value = 3011
$ws.Cells.Item(1106, 4).Value = 'Reko Swish +46732447361'
$ws.Cells.Item(1106, 6).Value = 495.54

# Row 1107
$ws.Cells.Item(1107, 1).Value = 44347
$ws.Cells.Item(1107, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1107, 2).Value = 'Reko194'
$ws.Cells.Item(1107, 3).Value = 2611
$ws.Cells.Item(1107, 4).Value = 'Reko Swish +46732447361'
$ws.Cells.Item(1107, 6).Value = 59.46

# Row 1108
$ws.Cells.Item(1108, 1).Value = 44347
$ws.Cells.Item(1108, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1108, 2).Value = 'Reko194'
$ws.Cells.Item(1108, 3).Value = 1930
$ws.Cells.Item(1108, 4).Value = 'Reko Swish +46732447361'
$ws.Cells.Item(1108, 5).Value = 555

# Row 1109
$ws.Cells.Item(1109, 1).Value = 44347
$ws.Cells.Item(1109, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1109, 2).Value = 'Reko195'
$ws.Cells.Item(1109, 3).Value = 3011
$ws.Cells.Item(1109, 4).Value = 'Reko Swish +46702179776'
$ws.Cells.Item(1109, 6).Value = 230.36

# Row 1110
$ws.Cells.Item(1110, 1).Value = 44347
$ws.Cells.Item(1110, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1110, 2).Value = 'Reko195'
$ws.Cells.Item(1110, 3).Value = 2611
$ws.Cells.Item(1110, 4).Value = 'Reko Swish +46702179776'
$ws.Cells.Item(1110, 6).Value = 27.64

# Row 1111
$ws.Cells.Item(1111, 1).Value = 44347
$ws.Cells.Item(1111, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1111, 2).Value = 'Reko195'
$ws.Cells.Item(1111, 3).Value = 1930
$ws.Cells.Item(1111, 4).Value = 'Reko Swish +46702179776'
$ws.Cells.Item(1111, 5).Value = 258

# Row 1112
$ws.Cells.Item(1112, 1).Value = 44347
$ws.Cells.Item(1112, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1112, 2).Value = 'Reko196'
$ws.Cells.Item(1112, 3).Value = 3011
$ws.Cells.Item(1112, 4).Value = 'Reko Swish +46707473778'
$ws.Cells.Item(1112, 6).Value = 423.21

# Row 1113
$ws.Cells.Item(1113, 1).Value = 44347
$ws.Cells.Item(1113, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1113, 2).Value = 'Reko196'
$ws.Cells.Item(1113, 3).Value = 2611
$ws.Cells.Item(1113, 4).Value = 'Reko Swish +46707473778'
$ws.Cells.Item(1113, 6).Value = 50.79

# Row 1114
$ws.Cells.Item(1114, 1).Value = 44347
$ws.Cells.Item(1114, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1114, 2).Value = 'Reko196'
$ws.Cells.Item(1114, 3).Value = 1930
$ws.Cells.Item(1114, 4).Value = 'Reko Swish +46707473778'
$ws.Cells.Item(1114, 5).Value = 474

# Row 1115
$ws.Cells.Item(1115, 1).Value = 44348
$ws.Cells.Item(1115, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1115, 2).Value = 'Reko197'
$ws.Cells.Item(1115, 3).Value = 3011
$ws.Cells.Item(1115, 4).Value = 'Reko Swish +46731587473'
$ws.Cells.Item(1115, 6).Value = 230.36

# Row 1116
$ws.Cells.Item(1116, 1).Value = 44348
$ws.Cells.Item(1116, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1116, 2).Value = 'Reko197'
$ws.Cells.Item(1116, 3).Value = 2611
$ws.Cells.Item(1116, 4).Value = 'Reko Swish +46731587473'
$ws.Cells.Item(1116, 6).Value = 27.64

# Row 1117
$ws.Cells.Item(1117, 1).Value = 44348
$ws.Cells.Item(1117, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1117, 2).Value = 'Reko197'
$ws.Cells.Item(1117, 3).Value = 1930
$ws.Cells.Item(1117, 4).Value = 'Reko Swish +46731587473'
$ws.Cells.Item(1117, 5).Value = 258

# Row 1118
$ws.Cells.Item(1118, 1).Value = 44348
$ws.Cells.Item(1118, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1118, 2).Value = 'Reko198'
$ws.Cells.Item(1118, 3).Value = 3011
$ws.Cells.Item(1118, 4).Value = 'Reko Swish +46703677212'
$ws.Cells.Item(1118, 6).Value = 460.71

# Row 1119
$ws.Cells.Item(1119, 1).Value = 44348
$ws.Cells.Item(1119, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1119, 2).Value = 'Reko198'
$ws.Cells.Item(1119, 3).Value = 2611
$ws.Cells.Item(1119, 4).Value = 'Reko Swish +46703677212'
$ws.Cells.Item(1119, 6).Value = 55.29

# Row 1120
$ws.Cells.Item(1120, 1).Value = 44348
$ws.Cells.Item(1120, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1120, 2).Value = 'Reko198'
$ws.Cells.Item(1120, 3).Value = 1930
$ws.Cells.Item(1120, 4).Value = 'Reko Swish +46703677212'
$ws.Cells.Item(1120, 5).Value = 516

# Row 1121
$ws.Cells.Item(1121, 1).Value = 44348
$ws.Cells.Item(1121, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1121, 2).Value = 'Reko199'
$ws.Cells.Item(1121, 3).Value = 3011
$ws.Cells.Item(1121, 4).Value = 'Reko Swish +46739881331'
$ws.Cells.Item(1121, 6).Value = 345.54

# Row 1122
$ws.Cells.Item(1122, 1).Value = 44348
$ws.Cells.Item(1122, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1122, 2).Value = 'Reko199'
$ws.Cells.Item(1122, 3).Value = 2611
$ws.Cells.Item(1122, 4).Value = 'Reko Swish +46739881331'
$ws.Cells.Item(1122, 6).Value = 41.46

# Row 1123
$ws.Cells.Item(1123, 1).Value = 44348
$ws.Cells.Item(1123, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1123, 2).Value = 'Reko199'
$ws.Cells.Item(1123, 3).Value = 1930
$ws.Cells.Item(1123, 4).Value = 'Reko Swish +46739881331'
$ws.Cells.Item(1123, 5).Value = 387

# Row 1124
$ws.Cells.Item(1124, 1).Value = 44349
$ws.Cells.Item(1124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1124, 2).Value = 'Reko200'
$ws.Cells.Item(1124, 3).Value = 3011
$ws.Cells.Item(1124, 4).Value = 'Reko Swish +46769256976'
$ws.Cells.Item(1124, 6).Value = 70.54000000000001

# Row 1125
$ws.Cells.Item(1125, 1).Value = 44349
$ws.Cells.Item(1125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1125, 2).Value = 'Reko200'
$ws.Cells.Item(1125, 3).Value = 2611
$ws.Cells.Item(1125, 4).Value = 'Reko Swish +46769256976'
$ws.Cells.Item(1125, 6).Value = 8.460000000000001

# Row 1126
$ws.Cells.Item(1126, 1).Value = 44349
$ws.Cells.Item(1126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1126, 2).Value = 'Reko200'
$ws.Cells.Item(1126, 3).Value = 1930
$ws.Cells.Item(1126, 4).Value = 'Reko Swish +46769256976'
$ws.Cells.Item(1126, 5).Value = 79

# Row 1127
$ws.Cells.Item(1127, 1).Value = 44349
$ws.Cells.Item(1127, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1127, 2).Value = 'Reko201'
$ws.Cells.Item(1127, 3).Value = 4010
$ws.Cells.Item(1127, 4).Value = 'Reko Swish +46731587473 return'
$ws.Cells.Item(1127, 5).Value = 115.18

# Row 1128
$ws.Cells.Item(1128, 1).Value = 44349
$ws.Cells.Item(1128, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1128, 2).Value = 'Reko201'
$ws.Cells.Item(1128, 3).Value = 2645
$ws.Cells.Item(1128, 4).Value = 'Reko Swish +46731587473 return'
$ws.Cells.Item(1128, 5).Value = 13.82

# Row 1129
$ws.Cells.Item(1129, 1).Value = 44349
$ws.Cells.Item(1129, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1129, 2).Value = 'Reko201'
$ws.Cells.Item(1129, 3).Value = 1930
$ws.Cells.Item(1129, 4).Value = 'Reko Swish +46731587473 return'
$ws.Cells.Item(1129, 6).Value = 129

# Row 1130
$ws.Cells.Item(1130, 1).Value = 44349
$ws.Cells.Item(1130, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1130, 2).Value = 'Reko202'
$ws.Cells.Item(1130, 3).Value = 3011
$ws.Cells.Item(1130, 4).Value = 'Reko Swish +46735426728'
$ws.Cells.Item(1130, 6).Value = 230.36

# Row 1131
$ws.Cells.Item(1131, 1).Value = 44349
$ws.Cells.Item(1131, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1131, 2).Value = 'Reko202'
$ws.Cells.Item(1131, 3).Value = 2611
$ws.Cells.Item(1131, 4).Value = 'Reko Swish +46735426728'
$ws.Cells.Item(1131, 6).Value = 27.64

# Row 1132
$ws.Cells.Item(1132, 1).Value = 44349
$ws.Cells.Item(1132, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1132, 2).Value = 'Reko202'
$ws.Cells.Item(1132, 3).Value = 1930
$ws.Cells.Item(1132, 4).Value = 'Reko Swish +46735426728'
$ws.Cells.Item(1132, 5).Value = 258

# Row 1133
$ws.Cells.Item(1133, 1).Value = 44349
$ws.Cells.Item(1133, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1133, 3).Value = 5670
$ws.Cells.Item(1133, 4).Value = 'ST1 V#LLINGBY K0135'
$ws.Cells.Item(1133, 5).Value = 757.17

# Row 1134
$ws.Cells.Item(1134, 1).Value = 44349
$ws.Cells.Item(1134, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1134, 3).Value = 2641
$ws.Cells.Item(1134, 4).Value = 'ST1 V#LLINGBY K0135'
$ws.Cells.Item(1134, 5).Value = 189.29

# Row 1135
$ws.Cells.Item(1135, 1).Value = 44349
$ws.Cells.Item(1135, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1135, 3).Value = 1930
$ws.Cells.Item(1135, 4).Value = 'ST1 V#LLINGBY K0135'
$ws.Cells.Item(1135, 6).Value = 946.46

# Row 1136
$ws.Cells.Item(1136, 1).Value = 44350
$ws.Cells.Item(1136, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1136, 2).Value = "'3030924"
$ws.Cells.Item(1136, 2).Style = "Normal"
$ws.Cells.Item(1136, 3).Value = 3011
$ws.Cells.Item(1136, 4).Value = 'Order 3030924 Swish +46703564388'
$ws.Cells.Item(1136, 6).Value = 1062.5

# Row 1137
$ws.Cells.Item(1137, 1).Value = 44350
$ws.Cells.Item(1137, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1137, 2).Value = "'3030924"
$ws.Cells.Item(1137, 2).Style = "Normal"
$ws.Cells.Item(1137, 3).Value = 2611
$ws.Cells.Item(1137, 4).Value = 'Order 3030924 Swish +46703564388'
$ws.Cells.Item(1137, 6).Value = 127.5

# Row 1138
$ws.Cells.Item(1138, 1).Value = 44350
$ws.Cells.Item(1138, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1138, 2).Value = "'3030924"
$ws.Cells.Item(1138, 2).Style = "Normal"
$ws.Cells.Item(1138, 3).Value = 1930
$ws.Cells.Item(1138, 4).Value = 'Order 3030924 Swish +46703564388'
$ws.Cells.Item(1138, 5).Value = 1190

# Row 1139
$ws.Cells.Item(1139, 1).Value = 44350
$ws.Cells.Item(1139, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1139, 2).Value = 'Reko203'
$ws.Cells.Item(1139, 3).Value = 3011
$ws.Cells.Item(1139, 4).Value = 'Reko Swish +46701825067'
$ws.Cells.Item(1139, 6).Value = 642.86

# Row 1140
$ws.Cells.Item(1140, 1).Value = 44350
$ws.Cells.Item(1140, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1140, 2).Value = 'Reko203'
$ws.Cells.Item(1140, 3).Value = 3011
$ws.Cells.Item(1140, 4).Value = 'Reko Swish +46701825067'
$ws.Cells.Item(1140, 6).Value = 8.039999999999999

# Row 1141
$ws.Cells.Item(1141, 1).Value = 44350
$ws.Cells.Item(1141, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1141, 2).Value = 'Reko203'
$ws.Cells.Item(1141, 3).Value = 2611
$ws.Cells.Item(1141, 4).Value = 'Reko Swish +46701825067'
$ws.Cells.Item(1141, 6).Value = 77.14

# Row 1142
$ws.Cells.Item(1142, 1).Value = 44350
$ws.Cells.Item(1142, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1142, 2).Value = 'Reko204'
$ws.Cells.Item(1142, 3).Value = 2611
$ws.Cells.Item(1142, 4).Value = 'Reko Swish +46701825067'
$ws.Cells.Item(1142, 6).Value = 0.96

# Row 1143
$ws.Cells.Item(1143, 1).Value = 44350
$ws.Cells.Item(1143, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1143, 2).Value = 'Reko204'
$ws.Cells.Item(1143, 3).Value = 1930
$ws.Cells.Item(1143, 4).Value = 'Reko Swish +46701825067'
$ws.Cells.Item(1143, 5).Value = 720

# Row 1144
$ws.Cells.Item(1144, 1).Value = 44350
$ws.Cells.Item(1144, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1144, 2).Value = 'Reko204'
$ws.Cells.Item(1144, 3).Value = 1930
$ws.Cells.Item(1144, 4).Value = 'Reko Swish +46701825067'
$ws.Cells.Item(1144, 5).Value = 9

# Row 1145
$ws.Cells.Item(1145, 1).Value = 44350
$ws.Cells.Item(1145, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1145, 2).Value = 'Reko205'
$ws.Cells.Item(1145, 3).Value = 3011
$ws.Cells.Item(1145, 4).Value = 'Reko Swish +46706395839'
$ws.Cells.Item(1145, 6).Value = 115.18

# Row 1146
$ws.Cells.Item(1146, 1).Value = 44350
$ws.Cells.Item(1146, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1146, 2).Value = 'Reko205'
$ws.Cells.Item(1146, 3).Value = 2611
$ws.Cells.Item(1146, 4).Value = 'Reko Swish +46706395839'
$ws.Cells.Item(1146, 6).Value = 13.82

# Row 1147
$ws.Cells.Item(1147, 1).Value = 44350
$ws.Cells.Item(1147, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1147, 2).Value = 'Reko205'
$ws.Cells.Item(1147, 3).Value = 1930
$ws.Cells.Item(1147, 4).Value = 'Reko Swish +46706395839'
$ws.Cells.Item(1147, 5).Value = 129

# Row 1148
$ws.Cells.Item(1148, 1).Value = 44350
$ws.Cells.Item(1148, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1148, 3).Value = 6400
$ws.Cells.Item(1148, 4).Value = 'FACEBK 62FYW4KZ62 K6885'
$ws.Cells.Item(1148, 5).Value = 415

# Row 1149
$ws.Cells.Item(1149, 1).Value = 44350
$ws.Cells.Item(1149, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1149, 4).Value = 'FACEBK 62FYW4KZ62 K6885'
$ws.Cells.Item(1149, 5).Value = 0

# Row 1150
$ws.Cells.Item(1150, 1).Value = 44350
$ws.Cells.Item(1150, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1150, 3).Value = 1930
$ws.Cells.Item(1150, 4).Value = 'FACEBK 62FYW4KZ62 K6885'
$ws.Cells.Item(1150, 6).Value = 415

# Row 1151
$ws.Cells.Item(1151, 1).Value = 44351
$ws.Cells.Item(1151, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1151, 3).Value = 6540
$ws.Cells.Item(1151, 4).Value = 'TELIA K0135'
$ws.Cells.Item(1151, 5).Value = 79.2

# Row 1152
$ws.Cells.Item(1152, 1).Value = 44351
$ws.Cells.Item(1152, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1152, 3).Value = 2641
$ws.Cells.Item(1152, 4).Value = 'TELIA K0135'
$ws.Cells.Item(1152, 5).Value = 19.8

# Row 1153
$ws.Cells.Item(1153, 1).Value = 44351
$ws.Cells.Item(1153, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1153, 3).Value = 1930
$ws.Cells.Item(1153, 4).Value = 'TELIA K0135'
$ws.Cells.Item(1153, 6).Value = 99

# Row 1154
$ws.Cells.Item(1154, 1).Value = 44352
$ws.Cells.Item(1154, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1154, 3).Value = 1220
$ws.Cells.Item(1154, 4).Value = 'CDON SE K0135'
$ws.Cells.Item(1154, 5).Value = 10463.2

# Row 1155
$ws.Cells.Item(1155, 1).Value = 44352
$ws.Cells.Item(1155, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1155, 3).Value = 2641
$ws.Cells.Item(1155, 4).Value = 'CDON SE K0135'
$ws.Cells.Item(1155, 5).Value = 2615.8

# Row 1156
$ws.Cells.Item(1156, 1).Value = 44352
$ws.Cells.Item(1156, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1156, 3).Value = 1930
$ws.Cells.Item(1156, 4).Value = 'CDON SE K0135'
$ws.Cells.Item(1156, 6).Value = 13079
